# Apply the edit described by the diff:
#  - The "Run 50" column (AZ) is removed entirely (header + data).
#    This shifts the old "Mean" column (BA) left into AZ.
#  - The header in A1 changes from "Gen" to "MaxFES".
#  - Column A values (generations) become fractional MaxFES values.
#  - The (now-shifted) Mean column values are updated to the recalculated
#    means (averages over the remaining 50 runs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the "Run 50" column (AZ). This shifts "Mean" (BA) into AZ.
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 2. Update the header for column A.
$ws.Range("A1").Value = "MaxFES"

# 3. Update column A data values (Gen -> MaxFES fractions).
$colA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
}

# 4. Update the (now shifted) Mean column values in AZ to the new recalculated means.
$means = @(138.78213951, 122.23217318, 66.11707583, 19.19769616, 11.1157757, 7.77381551, 6.1314417, 5.13540224, 4.18598638, 3.45876884, 3.04949454, 2.6642267, 2.47783992)
for ($i = 0; $i -lt $means.Length; $i++) {
    $row = $i + 2
    $ws.Range("AZ" + $row).Value = $means[$i]
}
